$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "stopStatement" (row 45) is no longer fully done -> completion degree drops
# from 1 (100%) to 0.5 (50%). B1's SUM(B2:B50)/49 average recalculates
# automatically.
$ws.Range("B45").Value = 0.5

# Because the statement is no longer fully completed, its row no longer gets
# the "done" (green) highlight; it now matches the "in progress" (yellow)
# look already used by rows such as A2. Copy that cell's formatting onto A45.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A45").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# Keep the active cell / view on the edited row.
$ws.Range("A45").Select()
